$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A21").Value = '$excuteCustomMethod'
$ws.Range("B21").Value = "Execute Custom Method"

$ws.Range("A21").Select()
